$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.949.88'
$ws.Cells.Item(2, 5).Value = '  +5.62%  '
$ws.Cells.Item(3, 4).Value = '2.256.17'
$ws.Cells.Item(3, 5).Value = '  +1.92%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'302.11"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.92%  '
$cell = $ws.Cells.Item(6, 4)
$cell.Value = "'92.36"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +6.45%  '
$cell = $ws.Cells.Item(7, 4)
$cell.Value = "'0.531"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +3.62%  '
$ws.Cells.Item(8, 5).Value = '  -0.06%  '
$ws.Cells.Item(9, 5).Value = '  +3.56%  '
$cell = $ws.Cells.Item(10, 4)
$cell.Value = "'54.54"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +9.28%  '
$cell = $ws.Cells.Item(11, 4)
$cell.Value = "'32.53"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +7.33%  '
$ws.Cells.Item(12, 5).Value = '  +2.42%  '
$ws.Cells.Item(13, 5).Value = '  +3.15%  '
$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'6.68"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +3.72%  '
$ws.Cells.Item(15, 4).Value = '2.602.60'
$ws.Cells.Item(15, 5).Value = '  +1.52%  '
$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'14.13"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +2.82%  '
$ws.Cells.Item(17, 4).Value = '2.274.49'
$ws.Cells.Item(17, 5).Value = '  +0.78%  '
$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'0.756"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +3.81%  '
$ws.Cells.Item(19, 4).Value = '41.849.07'
$ws.Cells.Item(19, 5).Value = '  +5.35%  '
$cell = $ws.Cells.Item(20, 4)
$cell.Value = "'12.15"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +9.95%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0901'
$ws.Cells.Item(21, 5).Value = '  +1.94%  '
$cell = $ws.Cells.Item(22, 4)
$cell.Value = "'5.93"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +3.69%  '
$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'67.09"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.35%  '
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'240.88"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.80%  '
$ws.Cells.Item(25, 5).Value = '  +5.51%  '
$ws.Cells.Item(26, 5).Value = '  -0.03%  '
$cell = $ws.Cells.Item(27, 4)
$cell.Value = "'1.90"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +4.23%  '
$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'23.89"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +4.11%  '
$ws.Cells.Item(29, 5).Value = '  +7.44%  '
$cell = $ws.Cells.Item(30, 4)
$cell.Value = "'9.66"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +4.92%  '
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'159.34"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.82%  '
$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'33.90"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +6.71%  '
$cell = $ws.Cells.Item(33, 4)
$cell.Value = "'0.999"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.03%  '
$ws.Cells.Item(34, 5).Value = '  +4.17%  '
$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'0.0743"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +4.60%  '
$cell = $ws.Cells.Item(36, 4)
$cell.Value = "'3.04"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +2.86%  '
$ws.Cells.Item(37, 5).Value = '  +2.89%  '
$ws.Cells.Item(38, 5).Value = '  +5.52%  '
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'16.51"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +8.11%  '
$ws.Cells.Item(40, 5).Value = '  +3.91%  '
$ws.Cells.Item(41, 5).Value = '  +3.50%  '
$ws.Cells.Item(42, 5).Value = '  +5.62%  '
$ws.Cells.Item(43, 4).Value = '2.048.24'
$ws.Cells.Item(43, 5).Value = '  -2.76%  '
$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'19.78"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +9.11%  '
$ws.Cells.Item(45, 5).Value = '  +3.46%  '
$ws.Cells.Item(46, 5).Value = '  +2.08%  '
$ws.Cells.Item(47, 5).Value = '  +2.32%  '
$cell = $ws.Cells.Item(48, 4)
$cell.Value = "'2.85"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +5.65%  '
$cell = $ws.Cells.Item(49, 4)
$cell.Value = "'1.52"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.79%  '
$ws.Cells.Item(50, 5).Value = '  +4.03%  '
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'51.81"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +6.39%  '
